$wb = $excel.ActiveWorkbook

# Rename the "linearized" sheet to "input_data"
$wsInput = $wb.Worksheets.Item("linearized")
$wsInput.Name = "input_data"

# The previously active/selected sheet was "art_initiation_rate" (tab 0);
# selection there stays at K5 (unchanged).
$wsRate = $wb.Worksheets.Item("art_initiation_rate")
$wsRate.Range("K5").Select()

# Active tab moves to the renamed "input_data" sheet (tab index 1), and its
# selection moves from D13:D17 to E25.
$wsInput.Activate()
$wsInput.Range("E25").Select()
